# Auto-generated script applying F-column ("想去人数") value updates
# as described by the commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 201
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 769
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 299
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 47
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 201
$ws.Cells.Item(3, 6).Value = 1375
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(19, 6).Value = 345
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 299
$ws.Cells.Item(29, 6).Value = 155
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(31, 6).Value = 549
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 12413
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 0

Write-Host "Updated column F values across sheets: 展览, 演出, 全部类型"
